$d = $word.ActiveDocument

# Locate the end of the existing "Whittaker (1962)" sentence.
$rng = $d.Content
$rng.Find.Execute("Whittaker (1962)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Text that leads up to the cross-reference to Figure 3 (fig:biomes).
$pre = " suggested that natural communities can be partitioned across biomes, largely defined as a function of their relative precipitation and temperature; in "
$rng.InsertAfter($pre)
$rng.Collapse(0)

# Insert a placeholder character, select it, then turn it into a hyperlink
# that cross-references the "fig:biomes" bookmark and displays "3".
$rng.InsertAfter("3")
$linkRng = $d.Range($rng.End - 1, $rng.End)
$d.Hyperlinks.Add($linkRng, "", "fig:biomes", "", "3")

# Continue after the hyperlink with the remainder of the paragraph.
$rng2 = $d.Range($linkRng.End, $linkRng.End)
$post = ", we show that even though networks, overall, capture the diversity of the precipitation/temperature climate well, types of networks have been studied in sub-spaces only. Specifically, parasitism networks have been studied in colder and drier climates; mutualism networks in wetter climates; predation networks display less of a bias."
$rng2.InsertAfter($post)

Write-Output "done"
